$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
}

Set-TextValue "D2" "243.89"
Set-TextValue "E2" "-0.07%"
Set-TextValue "D3" "29.94"
Set-TextValue "E3" "13.35%"
Set-TextValue "D4" "5.161"
Set-TextValue "E4" "0.29%"
Set-TextValue "E5" "1.22%"
Set-TextValue "D6" "6.542"
Set-TextValue "D7" "0.8480"
Set-TextValue "E7" "3.52%"
Set-TextValue "D8" "0.8661"
Set-TextValue "E8" "4.64%"
$ws.Range("B9").Value = "WazirX"
$ws.Range("C9").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextValue "D9" "0.1347"
Set-TextValue "E9" "1.25%"
$ws.Range("B10").Value = "MandalaExchangeToken"
$ws.Range("C10").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue "D10" "0.06909"
Set-TextValue "E10" "-0.09%"
$ws.Range("B11").Value = "BitrueCoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue "D11" "0.02918"
Set-TextValue "E11" "0.73%"
$ws.Range("B12").Value = "BitMartToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue "D12" "0.09374"
Set-TextValue "E12" "-0.13%"
$ws.Range("B13").Value = "BitForexToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue "D13" "0.001511"
$ws.Range("B14").Value = "CoinExToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
Set-TextValue "D14" "0.04177"
Set-TextValue "E14" "-9.31%"
$ws.Range("B15").Value = "One"
$ws.Range("C15").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextValue "D15" "0.0006009"
Set-TextValue "E15" "0.34%"
Set-TextValue "D16" "0.006034"
Set-TextValue "E16" "-3.17%"
Set-TextValue "D17" "3.511"
Set-TextValue "E17" "-3.95%"
Set-TextValue "E18" "0.00%"
Set-TextValue "D19" "2.245"
Set-TextValue "E19" "9.10%"
Set-TextValue "E21" "9.37%"
Set-TextValue "E22" "0.29%"
Set-TextValue "D23" "3.626"
Set-TextValue "E23" "-3.33%"
Set-TextValue "E25" "-1.17%"
Set-TextValue "D26" "0.004442"
Set-TextValue "E26" "-1.03%"
Set-TextValue "E27" "22.92%"
Set-TextValue "E28" "-0.56%"
Set-TextValue "D40" "0.03785"
Set-TextValue "E40" "3.95%"
Set-TextValue "D41" "0.005858"
Set-TextValue "E41" "72.58%"
Set-TextValue "D42" "0.1056"
Set-TextValue "E42" "-22.42%"
Set-TextValue "E43" "-11.27%"
Set-TextValue "D44" "0.009282"
Set-TextValue "E44" "2.84%"
Set-TextValue "D45" "0.00005101"
Set-TextValue "E45" "-4.90%"
Set-TextValue "E46" "0.00%"
Set-TextValue "E48" "-6.52%"
Set-TextValue "E49" "0.00%"
Set-TextValue "E50" "0.00%"
